$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for new column C
$ws.Range("C1").Value = "waterYear"

# Water year values for rows 2-12
$waterYears = @(2013,2014,2015,2016,2017,2018,2019,2020,2021,2022,2023)
for ($i = 0; $i -lt $waterYears.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $waterYears[$i]
}

# Update selection to match the diff (E14)
$ws.Range("E14").Select()
